# calorimetry : input and output consistency : done
#
# The free species (H, L, OH) rows are dropped from the output/consistency
# tables so only the fitted complexes (HL, H2L, HOH, HOHD) remain, matching
# the input tables. The enthalpies_calculated sheet also gets its sign
# convention corrected for the HL / H2L reaction enthalpies.

$wb = $excel.ActiveWorkbook

# --- input_enthalpies: drop H, L, OH rows -----------------------------
$wsEnthIn = $wb.Worksheets.Item("input_enthalpies")
$wsEnthIn.Range("A2:A4").EntireRow.Delete()

# --- constants_evaluated: drop H, L, OH rows ---------------------------
$wsConst = $wb.Worksheets.Item("constants_evaluated")
$wsConst.Range("A2:A4").EntireRow.Delete()

# --- enthalpies_calculated: drop H, L, OH rows, fix HL/H2L sign --------
$wsEnthCalc = $wb.Worksheets.Item("enthalpies_calculated")
$wsEnthCalc.Range("A2:A4").EntireRow.Delete()
$wsEnthCalc.Range("B2").Value = 13.0743647447613
$wsEnthCalc.Range("B3").Value = 17.5113480273213
